{"js": "const sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\nconst sec = sections.items[0];\nconst h2 = sec.getHeader(\"FirstPage\");\nreturn { ok: true };\n", "ps1": "$d = $word.ActiveDocument\n$sec = $d.Sections(1)\n$hdr = $sec.Headers(1)  # wdHeaderFooterPrimary = 1\n$hdr.Range.Text = \"Hello\"\n"}
